$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix company name in row 3 (ticker OB:ASC -> OB:ABG)
$ws.Range("B3").Value = "ABG Sundal Collier Holding ASA (OB:ABG)"

# Apply the same updated metrics to both data rows (2 and 3)
foreach ($row in 2, 3) {
    $ws.Range("D$row").Value = 0.0302
    $ws.Range("E$row").Value = 0.0103

    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = 29.4
    $ws.Range("L$row").Value = 0.1828358208955224
    $ws.Range("M$row").Value = 22.2
    $ws.Range("N$row").Value = 0.06713032960387059
    $ws.Range("O$row").Value = 0.7551020408163267
    $ws.Range("P$row").Value = 18.6
    $ws.Range("Q$row").Value = 0.05624433020864833
    $ws.Range("R$row").Value = 0.6326530612244898
    $ws.Range("S$row").Value = 3.600000000000001
    $ws.Range("T$row").Value = 0.1621621621621622
    $ws.Range("U$row").Value = 80.3
    $ws.Range("V$row").Value = 0.2428182642878742
    $ws.Range("W$row").Value = 0.4060773480662983
    $ws.Range("X$row").Value = 0.03009132918402517
    $ws.Range("Y$row").Value = 0.3759860188822731
    $ws.Range("Z$row").Value = 4.334231805929919
    $ws.Range("AA$row").Value = 0
    $ws.Range("AB$row").Value = 0.02907653416405395
    $ws.Range("AC$row").Value = -0.02907653416405395
    $ws.Range("AD$row").Value = 23.8
    $ws.Range("AE$row").Value = 0
    $ws.Range("AF$row").Value = 23.8
    $ws.Range("AG$row").Value = -56.5
    $ws.Range("AH$row").Value = 0.06713681241184767
    $ws.Range("AI$row").Value = 0.2224299065420561
    $ws.Range("AJ$row").Value = -0.2060539752005835
    $ws.Range("AK$row").Value = -2.116104868913857
    $ws.Range("AM$row").Value = -3.01

    # These columns are no longer populated for this row
    $ws.Range("AN$row").ClearContents()
    $ws.Range("AP$row").ClearContents()
}
